$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 0.5555555555555556
$ws.Range("D3").Value = 0.5555555555555556
$ws.Range("E3").Value = 0.8888888888888888
$ws.Range("F3").Value = 1
$ws.Range("H3").Value = 0.4159509202453988
$ws.Range("I3").Value = 0.2300363287408321
$ws.Range("J3").Value = 0.4444444444444444
$ws.Range("K3").Value = 377.8888888888889

$ws.Range("Q3").Value = 1063
$ws.Range("R3").Value = 20
$ws.Range("S3").Value = 137
$ws.Range("T3").Value = 526
$ws.Range("U3").Value = 673
$ws.Range("V3").Value = 558
$ws.Range("W3").Value = 1601
$ws.Range("X3").Value = 1484
$ws.Range("Y3").Value = 1095
$ws.Range("Z3").Value = 948

$ws.Range("AG3").Value = 0.987662
$ws.Range("AH3").Value = 0.915484
$ws.Range("AI3").Value = 0.675509
$ws.Range("AJ3").Value = 0.584824
